# ---------------------------------------------------------------------------
# Word COM-interop script implementing the commit:
#   "Added system model documentation closes #7"
#
# Changes applied:
#   1. The table-cell heading
#        "Jak ma wyglądać system? – 2 aplikacje na jednym serwerze + baza danych"
#      is rewritten to three separate (but identically formatted) runs and the
#      "2" is dropped while " + topic + queue" is appended, producing the
#      final text:
#        "Jak ma wyglądać system? –  aplikacje na jednym serwerze + baza danych + topic + queue"
#      split as: run1 "Jak ma wyglądać system? – "
#                run2 " aplikacje na jednym serwerze + baza danych"
#                run3 " + topic + queue"
#   2. The "Default Paragraph Font" style's uiPriority goes from 99 to 1 and
#      gains <w:unhideWhenUsed/>.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---- Part 1: split the heading run and append " + topic + queue" ----------
$searchText = "Jak ma wyglądać system? – 2 aplikacje na jednym serwerze + baza danych"
$searchRange = $d.Content
$found = $searchRange.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $searchRange.Start

    $run1 = "Jak ma wyglądać system? – "
    $run2 = " aplikacje na jednym serwerze + baza danych"
    $run3 = " + topic + queue"

    # Remember the original (shared) character formatting so the split can
    # restore it exactly once the run boundaries have been introduced.
    $origBold = $searchRange.Font.Bold

    # Boundary between run1 and the literal "2" that needs to be removed.
    $b1 = $start + $run1.Length
    $twoRange = $d.Range($b1, $b1 + 1)
    $twoRange.Delete()

    # Boundary between run2 and the new run3 text (after the "2" deletion).
    $b2 = $b1 + $run2.Length
    $tail = $d.Range($b2, $b2)
    $tail.InsertAfter($run3)

    $paraEnd = $b2 + $run3.Length

    # Force Word to materialize separate runs at each boundary by flipping a
    # character-formatting property and then flipping it straight back; the
    # resulting runs keep identical <w:rPr/> but remain split apart.
    $afterB1 = $d.Range($b1, $paraEnd)
    $afterB1.Font.Bold = -not $origBold
    $afterB1.Font.Bold = $origBold

    $afterB2 = $d.Range($b2, $paraEnd)
    $afterB2.Font.Bold = -not $origBold
    $afterB2.Font.Bold = $origBold

    $fullRange = $d.Range($start, $paraEnd)
    Write-Output "Heading updated: [$($fullRange.Text)]"
} else {
    Write-Output "WARNING: heading text was not found; no edit applied"
}

# ---- Part 2: "Default Paragraph Font" style tweaks ------------------------
$dpf = $d.Styles("Default Paragraph Font")
$dpf.Priority = 1
$dpf.UnhideWhenUsed = $true

Write-Output "done"
